{"js": "// Applies the Russian-translation edits described by the diff to Wallets.docx.\n// Each change is located with Body.search() (exact, case-sensitive text) and\n// replaced in place with Range.insertText(..., Word.InsertLocation.replace),\n// which keeps the run's original formatting (rPr) intact.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + searchText);\n  }\n\n  // Only the first occurrence in document order should be changed.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Web Wallet intro line (first occurrence only; the identical sentence in\n//    the Paper Wallet section further down stays in English).\nawait replaceOnce(\n  \"Please always check you are using the right URL\\u00A0\",\n  \"\u041f\u043e\u0436\u0430\u043b\u0443\u0439\u0441\u0442\u0430, \u0432\u0441\u0435\u0433\u0434\u0430 \u043f\u0440\u043e\u0432\u0435\u0440\u044f\u0439\u0442\u0435, \u0447\u0442\u043e \u0432\u044b \u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u0435\u0442\u0435 \u043f\u0440\u0430\u0432\u0438\u043b\u044c\u043d\u044b\u0439 URL \"\n);\n\n// 2) \"MOBILE WALLET\" heading.\nawait replaceOnce(\"MOBILE WALLET\", \"\u041c\u041e\u0411\u0418\u041b\u042c\u041d\u042b\u0415 \u041a\u041e\u0428\u0415\u041b\u042c\u041a\u0418\");\n\n// 3) Mobile wallet description line.\nawait replaceOnce(\n  \"SmartCash Mobile Wallets are run from your Mobile (Phone, Tablet,..) device\",\n  \"\u041c\u043e\u0431\u0438\u043b\u044c\u043d\u044b\u0435 \u043a\u043e\u0448\u0435\u043b\u044c\u043a\u0438 \u0434\u043b\u044f \u0432\u0430\u0448\u0435\u0433\u043e \u0443\u0441\u0442\u0440\u043e\u0439\u0441\u0442\u0432\u0430.\"\n);\n\n// 4) \"ELECTRUM WALLET\" heading.\nawait replaceOnce(\"ELECTRUM WALLET\", \"\u041a\u041e\u0428\u0415\u041b\u0401\u041a ELECTRUM\");\n\n// 5) Electrum wallet description line.\nawait replaceOnce(\n  \"This is a fast wallet that does not require the blockchain download. Wallet will not start SmartNodes, but a future release will add that feature.\",\n  \"\u0411\u044b\u0441\u0442\u0440\u044b\u0439 \u043a\u043e\u0448\u0435\u043b\u0451\u043a, \u043d\u0435 \u0442\u0440\u0435\u0431\u0443\u044e\u0449\u0438\u0439 \u0437\u0430\u0433\u0440\u0443\u0437\u043a\u0438 \u0431\u043b\u043e\u043a\u0447\u0435\u0439\u043d\u0430. \u0424\u0443\u043d\u043a\u0446\u0438\u044f \u0437\u0430\u043f\u0443\u0441\u043a\u0430 SmartNodes \u0431\u0443\u0434\u0435\u0442 \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d\u0430 \u0432 \u0431\u0443\u0434\u0443\u0449\u0438\u0445 \u0440\u0435\u043b\u0438\u0437\u0430\u0445.\"\n);\n\n// 6) \" SMARTCARD\" heading -> \"SMARTCARD\" (drop the leading NBSP).\nawait replaceOnce(\"\\u00A0SMARTCARD\", \"SMARTCARD\");\n\n// 7) SmartCard paragraph: only the first sentence is translated, the rest of\n//    the paragraph (already English) is left untouched.\nawait replaceOnce(\n  \"The SmartCard is a way to hold SMART on a physical card.\",\n  \"SmartCard \u2014 \u044d\u0442\u043e \u0441\u043f\u043e\u0441\u043e\u0431 \u0445\u0440\u0430\u043d\u0435\u043d\u0438\u044f SMART \u043d\u0430 \u0444\u0438\u0437\u0438\u0447\u0435\u0441\u043a\u043e\u0439 \u043a\u0430\u0440\u0442\u0435.\"\n);\n", "ps1": "# Applies the Russian-translation edits described by the diff to Wallets.docx.\n# Each change is located with Range.Find (exact, case-sensitive text, no\n# wildcards) starting from the top of the document; Find.Execute() collapses\n# $range to the matched text, and assigning $range.Text replaces it in place\n# while keeping the run's original formatting.\n#\n# NOTE: this host's script engine only reliably binds *positional* function\n# arguments, so Replace-FirstMatch takes $SearchText/$NewText positionally.\n\n$d = $word.ActiveDocument\n$nbsp = [char]0x00A0\n\nfunction Replace-FirstMatch($SearchText, $NewText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $SearchText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"No match found for: $SearchText\"\n    }\n\n    $range.Text = $NewText\n}\n\n# 1) Web Wallet intro line (first occurrence only; the identical sentence in\n#    the Paper Wallet section further down stays in English).\nReplace-FirstMatch (\"Please always check you are using the right URL\" + $nbsp) \"\u041f\u043e\u0436\u0430\u043b\u0443\u0439\u0441\u0442\u0430, \u0432\u0441\u0435\u0433\u0434\u0430 \u043f\u0440\u043e\u0432\u0435\u0440\u044f\u0439\u0442\u0435, \u0447\u0442\u043e \u0432\u044b \u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u0435\u0442\u0435 \u043f\u0440\u0430\u0432\u0438\u043b\u044c\u043d\u044b\u0439 URL \"\n\n# 2) \"MOBILE WALLET\" heading.\nReplace-FirstMatch \"MOBILE WALLET\" \"\u041c\u041e\u0411\u0418\u041b\u042c\u041d\u042b\u0415 \u041a\u041e\u0428\u0415\u041b\u042c\u041a\u0418\"\n\n# 3) Mobile wallet description line.\nReplace-FirstMatch \"SmartCash Mobile Wallets are run from your Mobile (Phone, Tablet,..) device\" \"\u041c\u043e\u0431\u0438\u043b\u044c\u043d\u044b\u0435 \u043a\u043e\u0448\u0435\u043b\u044c\u043a\u0438 \u0434\u043b\u044f \u0432\u0430\u0448\u0435\u0433\u043e \u0443\u0441\u0442\u0440\u043e\u0439\u0441\u0442\u0432\u0430.\"\n\n# 4) \"ELECTRUM WALLET\" heading.\nReplace-FirstMatch \"ELECTRUM WALLET\" \"\u041a\u041e\u0428\u0415\u041b\u0401\u041a ELECTRUM\"\n\n# 5) Electrum wallet description line.\nReplace-FirstMatch \"This is a fast wallet that does not require the blockchain download. Wallet will not start SmartNodes, but a future release will add that feature.\" \"\u0411\u044b\u0441\u0442\u0440\u044b\u0439 \u043a\u043e\u0448\u0435\u043b\u0451\u043a, \u043d\u0435 \u0442\u0440\u0435\u0431\u0443\u044e\u0449\u0438\u0439 \u0437\u0430\u0433\u0440\u0443\u0437\u043a\u0438 \u0431\u043b\u043e\u043a\u0447\u0435\u0439\u043d\u0430. \u0424\u0443\u043d\u043a\u0446\u0438\u044f \u0437\u0430\u043f\u0443\u0441\u043a\u0430 SmartNodes \u0431\u0443\u0434\u0435\u0442 \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d\u0430 \u0432 \u0431\u0443\u0434\u0443\u0449\u0438\u0445 \u0440\u0435\u043b\u0438\u0437\u0430\u0445.\"\n\n# 6) \" SMARTCARD\" heading -> \"SMARTCARD\" (drop the leading NBSP).\nReplace-FirstMatch ($nbsp + \"SMARTCARD\") \"SMARTCARD\"\n\n# 7) SmartCard paragraph: only the first sentence is translated, the rest of\n#    the paragraph (already English) is left untouched.\nReplace-FirstMatch \"The SmartCard is a way to hold SMART on a physical card.\" \"SmartCard \u2014 \u044d\u0442\u043e \u0441\u043f\u043e\u0441\u043e\u0431 \u0445\u0440\u0430\u043d\u0435\u043d\u0438\u044f SMART \u043d\u0430 \u0444\u0438\u0437\u0438\u0447\u0435\u0441\u043a\u043e\u0439 \u043a\u0430\u0440\u0442\u0435.\"\n"}
